# Updates cryptocurrency price (D) and volume-change (E) columns to match the
# latest scraped values from the GitHub Actions data-refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2 '26.264.98' -> '26.260.65'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.260.65'
$ws.Range("D2").Style = "Normal"
# Row 2: E2 '  -0.16%  ' -> '  -0.18%  '
$ws.Range("E2").Value = '  -0.18%  '

# Row 3: D3 '1.592.50' -> '1.591.51'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.591.51'
$ws.Range("D3").Style = "Normal"
# Row 3: E3 '  +0.16%  ' -> '  +0.12%  '
$ws.Range("E3").Value = '  +0.12%  '

# Row 4: E4 '  -0.04%  ' -> '  -0.02%  '
$ws.Range("E4").Value = '  -0.02%  '

# Row 5: D5 '212.47' -> '212.51'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.51'
$ws.Range("D5").Style = "Normal"
# Row 5: E5 '  +0.66%  ' -> '  +0.68%  '
$ws.Range("E5").Value = '  +0.68%  '

# Row 6: E6 '  -0.67%  ' -> '  -0.72%  '
$ws.Range("E6").Value = '  -0.72%  '

# Row 7: E7 '  -0.02%  ' -> '  +0.03%  '
$ws.Range("E7").Value = '  +0.03%  '

# Row 8: E8 '  -0.59%  ' -> '  -0.63%  '
$ws.Range("E8").Value = '  -0.63%  '

# Row 9: E9 '  -0.54%  ' -> '  -0.53%  '
$ws.Range("E9").Value = '  -0.53%  '

# Row 10: D10 '19.08' -> '19.04'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.04'
$ws.Range("D10").Style = "Normal"
# Row 10: E10 '  -1.69%  ' -> '  -1.81%  '
$ws.Range("E10").Value = '  -1.81%  '

# Row 11: D11 '0.0851' -> '0.0852'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0852'
$ws.Range("D11").Style = "Normal"
# Row 11: E11 '  +0.52%  ' -> '  +0.48%  '
$ws.Range("E11").Value = '  +0.48%  '

# Row 12: D12 '1.817.28' -> '1.815.50'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.815.50'
$ws.Range("D12").Style = "Normal"
# Row 12: E12 '  +0.14%  ' -> '  +0.11%  '
$ws.Range("E12").Value = '  +0.11%  '

# Row 13: D13 '1.596.46' -> '1.596.20'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.596.20'
$ws.Range("D13").Style = "Normal"
# Row 13: E13 '  +0.42%  ' -> '  +0.39%  '
$ws.Range("E13").Value = '  +0.39%  '

# Row 14: D14 '4.01' -> '4.00'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.00'
$ws.Range("D14").Style = "Normal"
# Row 14: E14 '  -1.77%  ' -> '  -1.67%  '
$ws.Range("E14").Value = '  -1.67%  '

# Row 15: E15 '  -2.40%  ' -> '  -2.44%  '
$ws.Range("E15").Value = '  -2.44%  '

# Row 16: D16 '63.88' -> '63.84'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.84'
$ws.Range("D16").Style = "Normal"

# Row 17: D17 '26.240.42' -> '26.247.75'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.247.75'
$ws.Range("D17").Style = "Normal"
# Row 17: E17 '  -0.25%  ' -> '  -0.24%  '
$ws.Range("E17").Value = '  -0.24%  '

# Row 18: E18 '  -0.63%  ' -> '  -0.64%  '
$ws.Range("E18").Value = '  -0.64%  '

# Row 19: D19 '216.05' -> '215.74'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.74'
$ws.Range("D19").Style = "Normal"
# Row 19: E19 '  +1.80%  ' -> '  +1.71%  '
$ws.Range("E19").Value = '  +1.71%  '

# Row 20: E20 '  -2.77%  ' -> '  -2.63%  '
$ws.Range("E20").Value = '  -2.63%  '

# Row 22: E22 '  +0.13%  ' -> '  +0.04%  '
$ws.Range("E22").Value = '  +0.04%  '

# Row 23: E23 '  +0.44%  ' -> '  +0.51%  '
$ws.Range("E23").Value = '  +0.51%  '

# Row 24: D24 '2.13' -> '2.12'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
# Row 24: E24 '  -1.10%  ' -> '  -1.23%  '
$ws.Range("E24").Value = '  -1.23%  '

# Row 25: D25 '144.35' -> '144.37'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.37'
$ws.Range("D25").Style = "Normal"
# Row 25: E25 '  +0.24%  ' -> '  +0.31%  '
$ws.Range("E25").Value = '  +0.31%  '

# Row 26: E26 '  -0.02%  ' -> '  +0.01%  '
$ws.Range("E26").Value = '  +0.01%  '

# Row 27: E27 '  -1.35%  ' -> '  -1.32%  '
$ws.Range("E27").Value = '  -1.32%  '

# Row 28: E28 '  -0.55%  ' -> '  -0.65%  '
$ws.Range("E28").Value = '  -0.65%  '

# Row 29: D29 '15.13' -> '15.12'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.12'
$ws.Range("D29").Style = "Normal"
# Row 29: E29 '  -0.73%  ' -> '  -0.91%  '
$ws.Range("E29").Value = '  -0.91%  '

# Row 30: D30 '0.0491' -> '0.0490'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0490'
$ws.Range("D30").Style = "Normal"
# Row 30: E30 '  -2.70%  ' -> '  -2.68%  '
$ws.Range("E30").Value = '  -2.68%  '

# Row 31: E31 '  +0.39%  ' -> '  +0.40%  '
$ws.Range("E31").Value = '  +0.40%  '

# Row 32: E32 '  -0.87%  ' -> '  -0.90%  '
$ws.Range("E32").Value = '  -0.90%  '

# Row 33: D33 '1.416.97' -> '1.415.16'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.415.16'
$ws.Range("D33").Style = "Normal"
# Row 33: E33 '  +6.76%  ' -> '  +6.52%  '
$ws.Range("E33").Value = '  +6.52%  '

# Row 35: E35 '  -0.38%  ' -> '  -0.36%  '
$ws.Range("E35").Value = '  -0.36%  '

# Row 36: E36 '  -0.70%  ' -> '  -0.95%  '
$ws.Range("E36").Value = '  -0.95%  '

# Row 37: D37 '0.584' -> '0.583'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.583'
$ws.Range("D37").Style = "Normal"
# Row 37: E37 '  -3.24%  ' -> '  -3.32%  '
$ws.Range("E37").Value = '  -3.32%  '

# Row 38: E38 '  -1.03%  ' -> '  -1.01%  '
$ws.Range("E38").Value = '  -1.01%  '

# Row 39: D39 '0.824' -> '0.823'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.823'
$ws.Range("D39").Style = "Normal"
# Row 39: E39 '  +0.77%  ' -> '  +0.74%  '
$ws.Range("E39").Value = '  +0.74%  '

# Row 40: D40 '5.88' -> '5.85'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.85'
$ws.Range("D40").Style = "Normal"
# Row 40: E40 '  +2.91%  ' -> '  +2.42%  '
$ws.Range("E40").Value = '  +2.42%  '

# Row 42: D42 '0.975' -> '0.974'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.974'
$ws.Range("D42").Style = "Normal"
# Row 42: E42 '  -2.05%  ' -> '  -1.74%  '
$ws.Range("E42").Value = '  -1.74%  '

# Row 43: E43 '  +0.17%  ' -> '  +0.19%  '
$ws.Range("E43").Value = '  +0.19%  '

# Row 44: D44 '0.765' -> '0.766'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.766'
$ws.Range("D44").Style = "Normal"
# Row 44: E44 '  -0.08%  ' -> '  +0.01%  '
$ws.Range("E44").Value = '  +0.01%  '

# Row 45: D45 '1.729.05' -> '1.728.08'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.728.08'
$ws.Range("D45").Style = "Normal"
# Row 45: E45 '  +0.08%  ' -> '  +0.11%  '
$ws.Range("E45").Value = '  +0.11%  '

# Row 46: D46 '60.99' -> '60.92'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.92'
$ws.Range("D46").Style = "Normal"
# Row 46: E46 '  -1.53%  ' -> '  -1.56%  '
$ws.Range("E46").Value = '  -1.56%  '

# Row 47: D47 '86.46' -> '86.21'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.21'
$ws.Range("D47").Style = "Normal"
# Row 47: E47 '  -1.90%  ' -> '  -2.15%  '
$ws.Range("E47").Value = '  -2.15%  '

# Row 48: D48 '1.49' -> '1.48'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.48'
$ws.Range("D48").Style = "Normal"
# Row 48: E48 '  +0.06%  ' -> '  +0.04%  '
$ws.Range("E48").Value = '  +0.04%  '

# Row 49: E49 '  -0.72%  ' -> '  -0.69%  '
$ws.Range("E49").Value = '  -0.69%  '

# Row 50: E50 '  -2.53%  ' -> '  -2.52%  '
$ws.Range("E50").Value = '  -2.52%  '

# Row 51: E51 '  -0.13%  ' -> '  -0.03%  '
$ws.Range("E51").Value = '  -0.03%  '
